# update_account_tags.xlsx - clear the sample/test row so the sheet is
# ready for the automation to populate real account data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the placeholder test data row (account number + shared-string
# Advisor/Payout names) while keeping the header row and A2's style intact.
$ws.Range("A2").ClearContents()
$ws.Range("B2:C2").ClearContents()

# Match the author's final selection state: the whole of row 2 selected.
[void]$ws.Rows("2:2").Select()
